$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data (rows 1-5):
#  1: Email | Name | Company | Sent                         (header)
#  2: verma.arpit078@gmail.com | Arpit | xyz | 23:57 13-06-2023
#  3: blabbla429@gmail.com | bla | google | 23:57 13-06-2023
#  4: 2021eeb1156@iitrpr.ac.in | college id | oracle | 23:57 13-06-2023
#  5: server1078@outlook.com | server | zerodha | 23:57 13-06-2023
#
# Target data (rows 1-3):
#  1: Email | Name | Company | Sent                         (header, unchanged)
#  2: blabbla429@gmail.com | bla | google | 0:43 15-06-2023
#  3: server1078@outlook.com | server | zerodha | 0:43 15-06-2023
#
# i.e. the "verma.arpit078" entry (row 2) and the "2021eeb1156" / college id
# entry (row 4) are removed, the remaining two rows shift up and get a
# refreshed "Sent" timestamp.

# Remove the verma.arpit078@gmail.com row (originally row 2). Remaining rows
# shift up by one.
$ws.Rows(2).Delete()

# The "2021eeb1156@iitrpr.ac.in" / college id row, originally row 4, is now
# row 3 after the previous delete. Remove it too.
$ws.Rows(3).Delete()

# Refresh the "Sent" timestamps on the two surviving rows.
$ws.Range("D2").Value = "0:43 15-06-2023"
$ws.Range("D3").Value = "0:43 15-06-2023"
